{"js": "// Streaming MJPEG and fixed refresh\n// Switched from JavaScript supported MJPEG to browser supported MJPEG,\n// which is about 3 times faster.  The raspicam invocation gains a\n// \"-hf\" (horizontal flip) flag right before the existing \"-fps 15 -q 50\"\n// options, i.e. \"-fps 15 -q 50\" becomes \"-hf -fps 15 -q 50\" (en dash\n// variant of the leading hyphen, matching the surrounding \"-vf\"/\"-ex\").\n\n// 1. Update the command-line text in the document body. Locate the\n//    unique \"-fps 15 -q 50 \" fragment (made up of the old \"-f\" run\n//    followed by the old \"ps 15 -q 50 \" run), then narrow down to the\n//    \"-f\" at its start and insert \"<en-dash>hf \" right before it. This\n//    turns \"-fps 15 -q 50 \" into \"<en-dash>hf -fps 15 -q 50 \" while\n//    leaving the existing runs (and the _GoBack bookmark sitting\n//    between them) untouched, instead of rewriting the whole match.\nconst target = context.document.body.search(\"-fps 15 -q 50 \", { matchCase: true });\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  const match = target.items[0];\n  const prefix = match.search(\"-f\", { matchCase: true });\n  prefix.load(\"items\");\n  await context.sync();\n\n  if (prefix.items.length > 0) {\n    prefix.items[0].insertText(\"\\u2013hf \", Word.InsertLocation.before);\n    await context.sync();\n  }\n}\n\n// 2. The inserted text pushes the document one line over a page\n//    boundary, so the cached PAGE field result shown in the footer\n//    advances from 6 to 7.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < sections.items.length; i++) {\n  const footer = sections.items[i].getFooter(Word.HeaderFooterType.primary);\n  const hits = footer.search(\"6\", { matchCase: true, matchWholeWord: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < hits.items.length; j++) {\n    hits.items[j].insertText(\"7\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Streaming MJPEG and fixed refresh\n# Switched from JavaScript supported MJPEG to browser supported MJPEG,\n# which is about 3 times faster.  The raspicam invocation gains a\n# \"-hf\" (horizontal flip) flag right before the existing \"-fps 15 -q 50\"\n# options, i.e. \"-fps 15 -q 50\" becomes \"-hf -fps 15 -q 50\" (en dash\n# variant of the leading hyphen, matching the surrounding \"-vf\"/\"-ex\").\n\n$d = $word.ActiveDocument\n\n# 1. Update the command-line text in the document body. Locate the\n#    unique \"-fps 15 -q 50 \" fragment (made up of the old \"-f\" run\n#    followed by the old \"ps 15 -q 50 \" run) and insert \"<en-dash>hf \"\n#    directly in front of it, turning \"-fps 15 -q 50 \" into\n#    \"<en-dash>hf -fps 15 -q 50 \". InsertBefore on the live Find range\n#    keeps the existing runs (and the _GoBack bookmark sitting between\n#    them) intact instead of rewriting the whole paragraph.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"-fps 15 -q 50 \"\n$found = $find.Execute()\nif ($found) {\n  $rng.InsertBefore([char]0x2013 + \"hf \")\n}\n\n# 2. The inserted text pushes the document one line over a page\n#    boundary, so the cached PAGE field result shown in the footer\n#    advances from 6 to 7.\nfor ($s = 1; $s -le $d.Sections.Count; $s++) {\n  $sec = $d.Sections.Item($s)\n  for ($i = 1; $i -le 3; $i++) {\n    $footer = $sec.Footers.Item($i)\n    foreach ($f in $footer.Range.Fields) {\n      if ($f.Code.Text.Trim() -eq \"PAGE\" -and $f.Result.Text -eq \"6\") {\n        $ffind = $footer.Range.Find\n        $ffind.ClearFormatting()\n        $ffind.Text = \"6\"\n        $ffind.Replacement.ClearFormatting()\n        $ffind.Replacement.Text = \"7\"\n        $ffind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1) | Out-Null\n      }\n    }\n  }\n}\n"}
